# Remove the "User phone/OTP" test-case rows that were added to Sheet1
# (H22:I27) and re-select the "User Registration_TestCases" tab.

$wb = $excel.ActiveWorkbook

$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Range("H22:I27").ClearContents()
$wsSheet1.Rows("22:27").AutoFit()

# Activate the "User Registration_TestCases" tab (this is what makes Excel
# flip tabSelected/activeTab/topLeftCell/selection bookkeeping between the
# two sheets, matching the diff).
$wsReg = $wb.Worksheets.Item("User Registration_TestCases")
$wsReg.Activate()
$wsReg.Range("H5").Select()

# Restore Sheet1's own scroll position / selection to match the target view
# (topLeftCell A21, selection H22:I27) without it being the active sheet.
$wsSheet1.Range("H22:I27").Select()
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1

$wsReg.Activate()
$wsReg.Range("H5").Select()
